$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused column B entirely (B1 header timestamp + empty B2)
$ws.Range("B1").EntireColumn.Delete()

# Append the new log entry in column A
$ws.Range("A3").Value = "2025-05-28 11:43:47"
